$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.0763005239001599
$ws.Range("B1").Value = -0.07630052526628002

$ws.Range("A2").Value = -0.063755276670105887
$ws.Range("B2").Value = 0.063755275294328922

$ws.Range("A3").Value = -0.035576518170188254
$ws.Range("B3").Value = 0.035576516778762793

$ws.Range("A4").Value = 0.05347974304410013
$ws.Range("B4").Value = -0.053479744490016336

$ws.Range("A5").Value = -0.025331793880731899
$ws.Range("B5").Value = 0.025331792427029917
